$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Productos")

# The single example row in "Productos" used "N/A" placeholders for
# SKU / Peso / Altura / Longitud / Profundidad (columns E-I) and the
# project shipped a second example workbook just to show real numbers
# for those fields (see "Variantes", row 2: 12345678 / 2000 / 200 / 10 /
# 10 / 10). Fill in realistic numeric sample data here instead, so one
# workbook is enough.
$ws.Range("E2").Value = 87654321
$ws.Range("F2").Value = 200
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 10

# Leave the selection on the last touched cell.
[void]$ws.Range("I2").Select()
